$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Updated) date column (C) for rows 2 through 79
# from serial date 45186 (2023-09-17) to serial date 45188 (2023-09-19).
for ($row = 2; $row -le 79; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45188
}
